# Insert a new data row at row 569 (pushing existing rows 569-610 down to 570-611)
# and populate it with the new reading: 2026/01/06, 火, 22, 166

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(569).Insert()

# Force column A to remain text so the date-like string "2026/01/06" is not
# auto-converted into a date serial number (matching the other rows in the
# sheet, which all store dates as literal text).
$ws.Range("A569").NumberFormat = "@"
$ws.Range("A569").Value = "2026/01/06"
$ws.Range("B569").Value = "火"
$ws.Range("C569").Value = 22
$ws.Range("D569").Value = 166
